$d = $word.ActiveDocument

# The third bullet under "Review of work undertaken" (bulleted list, numId 9)
# currently contains a single placeholder space character. Replace it with
# the new status text, keeping the run's existing formatting
# (rFonts/sz/szCs/u theme-minor formatting) untouched.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    # Range.Text includes the trailing paragraph-mark character, so strip
    # any trailing CR before comparing against the lone-space placeholder.
    $text = $cand.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq " ") {
        $target = $cand
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the placeholder paragraph (single space) to update."
}

$pRange = $target.Range

# Insert the new text immediately before the paragraph mark (i.e. right
# after the existing space), then delete the old leading space. Doing the
# insert first (rather than assigning straight over the space) keeps the
# new text merged into the existing run so its formatting/rPr is preserved
# and no spurious xml:space="preserve" is introduced.
$insertPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$insertPoint.InsertBefore("Foundations of live trading implemented.")

$updated = $target.Range
$oldSpace = $d.Range($updated.Start, $updated.Start + 1)
$oldSpace.Text = ""
